# Updates cryptos list values (Price and Volume(1h) columns) per the
# "Updated cryptos list" data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.205.96"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.922.22"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.94"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.65"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.500"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.39"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "3.406.69"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "61.206.41"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "2.925.17"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.66"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "431.35"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.50"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.06"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.70"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.89"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.73"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -4.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.59"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.89"
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.59"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "0.0₃0878"
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.62"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.99"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.122"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.51"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.280"
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0344"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("D44").Value = "2.692.22"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "133.55"
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "364.22"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.55"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("E51").Value = "  -0.67%  "
